# Scoreboard.xlsx update: add the 3rd-event (Minute3/Second3/Rep3 -> J/K/L)
# results for the "ScoreM" sheet, and refresh the active-sheet / selection
# state to match what was left after the edit (ScoreM active with L18
# selected, ScoreF's own selection moved to J10).

$wb  = $excel.ActiveWorkbook
$wsM = $wb.Worksheets.Item("ScoreM")
$wsF = $wb.Worksheets.Item("ScoreF")

# Newly-entered event-3 scores (Minute3=J, Second3=K, Rep3=L) for rows 2-23.
$event3 = @(
    @(2, 13, 0, 223),
    @(3, 13, 0, 226),
    @(4, 13, 0, 224),
    @(5, 11, 32, 258),
    @(6, 13, 0, 243),
    @(7, 13, 0, 176),
    @(8, 9, 49, 258),
    @(9, 13, 0, 213),
    @(10, 10, 47, 258),
    @(11, 13, 0, 240),
    @(12, 10, 14, 258),
    @(13, 13, 0, 203),
    @(14, 9, 55, 258),
    @(15, 10, 8, 258),
    @(16, 12, 9, 258),
    @(17, 11, 21, 258),
    @(18, 13, 0, 176),
    @(19, 11, 1, 258),
    @(20, 11, 31, 258),
    @(21, 10, 49, 258),
    @(22, 10, 23, 258),
    @(23, 13, 0, 237)
)

foreach ($entry in $event3) {
    $r = $entry[0]
    $wsM.Cells.Item($r, 10).Value = $entry[1]
    $wsM.Cells.Item($r, 11).Value = $entry[2]
    $wsM.Cells.Item($r, 12).Value = $entry[3]
}

# ScoreF keeps its own selection, just moved to J10, and loses the
# "tabSelected" flag as ScoreM becomes the active sheet.
$wsF.Activate()
$wsF.Range("J10").Select()

# ScoreM becomes the active sheet with L18 selected.
$wsM.Activate()
$wsM.Range("L18").Select()
